$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: classification of the adaptation/mitigation sentence ---
$ws.Range("C16").Value = "adaptation, mitigation, Millenium Development Goals"
$ws.Range("G16").Value = "sufficientarian"
$ws.Range("H16").Value = "explicitly referring to the reason for prescribing a balanced view of mitigation and adaptation to have the world reach sufficient targets. "

# --- Row 17: classification of the financial resources sentence ---
$ws.Range("F17").Value = "present"
$ws.Range("H17").Value = "prescribing the need of developing countries to fulfill funding of developing countries"

# --- Row 19: newly classified as relevant (women/gender sentence) ---
$ws.Range("B19").Value = "yes"
$ws.Range("C19").Value = "women, gender"
$ws.Range("D19").Value = "n.a."
$ws.Range("E19").Value = "global"
$ws.Range("F19").Value = "n.a."
$ws.Range("G19").Value = "egalitarian"
$ws.Range("H19").Value = "Contains the value judgement and prescription of including women in new policies. Presenting an egalitarian view. "

# --- Row 25: newly classified as relevant (Typhoon Haiyan / call to action sentence) ---
$ws.Range("B25").Value = "yes"
$ws.Range("C25").Value = "action"
$ws.Range("D25").Value = "n.a."
$ws.Range("E25").Value = "global"
$ws.Range("F25").Value = "present"
$ws.Range("G25").Value = "general normative statement"
$ws.Range("H25").Value = "Prescribing the need to take measures, no distinctive distribution preferred. "

# --- Update the active view: scroll/freeze anchor and current selection ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$ws.Range("C16").Select() | Out-Null
